$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
try {
  $excel.Goto($ws.Range("A12"), $true)
  Write-Host "Goto worked"
} catch {
  Write-Host ("Goto failed: " + $_)
}
$ws.Range("B14").Select()
